$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Title text box: "version: 30 May 2018: " -> split + date change ---
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$dateRange = $titleTr.Characters(32, 13)
$dateRange.Text = ": 1 June 2018"

# --- 2. "Sep/Okt 2018" -> "Sep/Oct 2018" (typo correction), split into runs ---
$sepShape = $s.Shapes.Item(4).GroupItems.Item(1)
$oktRange = $sepShape.TextFrame.TextRange.Characters(5, 3)
$oktRange.Text = "Oct"

# --- 3. Add new arrow/line shape (duplicate an existing similar shape) ---
$srcShape = $s.Shapes.Item(64)
$dupRange = $srcShape.Duplicate()
$newShape = $dupRange.Item(1)
$newShape.Left = 236.5840301513672
$newShape.Top = 1301.6431884765625
$newShape.Width = 2.0293703079223633
$newShape.Height = 110.43756103515625
